$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text type (prevents Excel from
# auto-converting numeric-looking strings like "550.68" into floating point
# numbers), then reset the style back to Normal so no stray style index is
# left referenced on the cell (keeps formatting identical to the original).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "57.446.11"
$ws.Range("E2").Value = "  -6.11%  "
$ws.Range("D3").Value = "2.897.81"
$ws.Range("E3").Value = "  -3.72%  "
$ws.Range("E4").Value = "  +0.27%  "
Set-TextValue "D5" "550.68"
$ws.Range("E5").Value = "  -2.38%  "
Set-TextValue "D6" "122.95"
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "2.886.40"
$ws.Range("E8").Value = "  -4.05%  "
Set-TextValue "D9" "0.496"
$ws.Range("E9").Value = "  -0.42%  "
Set-TextValue "D10" "0.124"
$ws.Range("E10").Value = "  -8.45%  "
Set-TextValue "D11" "4.70"
$ws.Range("E11").Value = "  -10.36%  "
Set-TextValue "D12" "0.438"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("E13").Value = "  -5.50%  "
Set-TextValue "D14" "32.44"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "3.374.08"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("D17").Value = "2.895.98"
$ws.Range("E17").Value = "  -3.52%  "
Set-TextValue "D18" "6.53"
$ws.Range("E18").Value = "  +4.42%  "
$ws.Range("D19").Value = "57.468.26"
$ws.Range("E19").Value = "  -5.97%  "
Set-TextValue "D20" "404.57"
$ws.Range("E20").Value = "  -7.78%  "
$ws.Range("E21").Value = "  -2.45%  "
Set-TextValue "D22" "0.671"
$ws.Range("E22").Value = "  +0.78%  "
Set-TextValue "D23" "6.83"
$ws.Range("E23").Value = "  -4.68%  "
Set-TextValue "D24" "12.75"
$ws.Range("E24").Value = "  +1.11%  "
Set-TextValue "D25" "77.01"
$ws.Range("E25").Value = "  -2.72%  "
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("E28").Value = "  -1.86%  "
Set-TextValue "D29" "1.94"
$ws.Range("E29").Value = "  +1.96%  "
Set-TextValue "D30" "7.18"
$ws.Range("E30").Value = "  -0.72%  "
Set-TextValue "D31" "6.03"
$ws.Range("E31").Value = "  -2.82%  "
Set-TextValue "D32" "24.65"
$ws.Range("E32").Value = "  -3.78%  "
Set-TextValue "D33" "0.0986"
$ws.Range("E33").Value = "  +4.75%  "
$ws.Range("E34").Value = "  -2.76%  "
Set-TextValue "D35" "0.907"
$ws.Range("E35").Value = "  -5.11%  "
$ws.Range("E36").Value = "  -12.29%  "
Set-TextValue "D37" "47.86"
$ws.Range("E37").Value = "  -4.63%  "
Set-TextValue "D38" "8.35"
$ws.Range("E38").Value = "  +7.44%  "
$ws.Range("D39").Value = "0.0₃0619"
$ws.Range("E39").Value = "  -7.96%  "
$ws.Range("E40").Value = "  -2.19%  "
Set-TextValue "D41" "0.0341"
$ws.Range("E41").Value = "  -5.80%  "
$ws.Range("D42").Value = "2.619.34"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.42"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D44" "360.34"
$ws.Range("E44").Value = "  -4.53%  "
Set-TextValue "D45" "0.998"
$ws.Range("E45").Value = "  -0.05%  "
Set-TextValue "D46" "118.53"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  -3.45%  "
Set-TextValue "D48" "0.106"
$ws.Range("E48").Value = "  -0.22%  "
Set-TextValue "D49" "1.95"
$ws.Range("E49").Value = "  -1.95%  "
Set-TextValue "D50" "22.85"
$ws.Range("E50").Value = "  -2.80%  "
$ws.Range("E51").Value = "  -4.67%  "
